$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new user-activity log rows (LOGIN events) below the
# existing data (rows 2-21), growing the used range to A1:D24.

$ws.Range("A22").Value = "2025-12-26 18:35:33"
$ws.Range("B22").Value = "manikandaa944@gmail.com"
$ws.Range("C22").Value = "LOGIN"
$ws.Range("D22").Value = "User logged in successfully"

$ws.Range("A23").Value = "2025-12-28 16:09:04"
$ws.Range("B23").Value = "manikandaa944@gmail.com"
$ws.Range("C23").Value = "LOGIN"
$ws.Range("D23").Value = "User logged in successfully"

$ws.Range("A24").Value = "2025-12-28 16:15:25"
$ws.Range("B24").Value = "manikandaa944@gmail.com"
$ws.Range("C24").Value = "LOGIN"
$ws.Range("D24").Value = "User logged in successfully"
